# Update cryptos list with latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.845.82"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "1.822.57"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3011"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07719"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("D12").Value = "1.827.20"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.029"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6701"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.413"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008252"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "28.840.57"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.60%  "
$ws.Range("D20").Value = "2.082.04"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.364"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1467"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.681"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.179"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.124"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.186"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05091"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7514"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.811"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.137"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01824"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").Value = "1.203.99"
$ws.Range("E39").Value = "  -4.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.674"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "1.982.02"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5158"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  -6.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.405"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.234"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -12.66%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.724"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "62.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4156"
$ws.Range("D51").Style = "Normal"
